$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68 holds the template record ("moses" / "bro" / "1234" / "m@g.c" /
# "Male" / 0) that this sheet already repeats many times above. Append
# eight more identical rows (69-76) below it, the same way the existing
# duplicate rows were produced - by copying the template row down.
$template = $ws.Range("A68:F68")

for ($row = 69; $row -le 76; $row++) {
    $target = $ws.Range("A" + $row + ":F" + $row)
    $template.Copy($target)
}
